$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

function Replace-Word($old, $new) {
    # whole-word match, to avoid accidentally hitting the string as a substring
    # of an unrelated word elsewhere in the document (e.g. "net" inside "planet")
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Title ---
Replace-Text "Symbiosis: Nature's Collaborative Harmony" "The Historical Significance of Arts and Culture: A Brief History Through the Ages"

# --- Author name ---
Replace-Text " Jessica Kennedy" " Mark Anthony"

# --- Email address (three runs) ---
Replace-Word "kennedy" "mark"
Replace-Text "jessica25@archetype" "anthony@eduschool"
Replace-Word "net" "org"

# --- Body paragraph 1, first sentence chain ---
Replace-Text "In the vast tapestry of life, symbiotic relationships paint a vibrant picture of interconnectedness and interdependence" "Art and culture serve as a mirror to society, reflecting the values, beliefs, and experiences of various civilizations throughout history"

Replace-Text "Like threads weaving together a complex web, diverse organisms engage in intricate partnerships, mutually benefiting from their shared existence" "From the cave paintings of early humans to the modern masterpieces that grace museums globally, art has captured the essence of human ingenuity and creativeness"

Replace-Text " It is within these cooperative arrangements that we find profound lessons on survival, cooperation, and the delicate balance that sustains our planet." ""

Replace-Text "From microscopic realms to expansive ecosystems, symbiosis manifests in myriad forms, offering insights into the intricate beauty of nature's collaborative harmony" "This essay explores the historical significance of arts and culture, highlighting their profound impact on shaping societies and fostering cultural identities across different eras"

# --- Body paragraph 1, second sentence chain ---
Replace-Text "Decoding the language of symbiosis reveals a symphony of interspecies interactions" "The earliest forms of art, such as cave paintings and sculptures, served as a means of communication and storytelling for ancient civilizations"

Replace-Text "Mutualism, a fundamental aspect of this biological ballet, involves the reciprocal exchange of benefits" "These artistic expressions provided insights into their daily lives, spiritual beliefs, and perceptions of the natural world"

Replace-Text "Like partners in a tango, organisms engage in mutually advantageous behaviors, enhancing each other's survival and flourishing" "As societies evolved, art became increasingly sophisticated, leading to the development of various art forms, including painting, sculpture, architecture, music, and literature"

Replace-Text " Commensalism, a less intimate form, exists when one organism benefits while the other remains unaffected, akin to a lodger sharing space rent-free. Parasitism, a more contentious alliance, sees one organism exploiting another for its own gain, a dynamic reminiscent of a predator and prey relationship." ""

Replace-Text "These symbiotic interactions, with their varying degrees of cooperation and exploitation, shape the very fabric of our ecosystems, contributing to biodiversity, stability, and the harmonious coexistence of life" "The Renaissance period witnessed an artistic explosion that celebrated humanism and individualism, producing timeless works of art that continue to inspire audiences today"

# --- Body paragraph 1, third sentence chain ---
Replace-Text "Venturing into the fascinating world of commensalism, we encounter organisms like the remora, a fish that attaches itself to larger marine creatures, enjoying protection and a free ride while posing no harm to its host" "In modern times, art and culture have become integral to the fabric of society, reflecting the multifaceted nature of the human experience"

Replace-Text "Epiphytes, plants that grow on the branches and trunks of trees without causing damage, exemplify this cooperative balance" " From the avant-garde movements of the 20th century to the rise of digital and interactive media today, art has embraced new technologies and mediums to push boundaries and explore new dimensions of creative expression"

Replace-Text " As they bask in their arboreal haven, epiphytes absorb nutrients from the air and rain, contributing to the overall productivity of the forest ecosystem. In a different realm, the human microbiome, a complex community of bacteria, fungi, and other microorganisms, forms a commensal partnership with its human host." ""

Replace-Text "These microorganisms play vital roles in digestion, immunity, and overall health, illustrating the intricate dance between host and microbe" "Contemporary art often engages with social, political, and environmental issues, challenging preconceived notions and provoking critical thought"

# --- Summary heading: lastRenderedPageBreak is a rendering cache artifact, left alone ---

# --- Summary paragraph ---
Replace-Text "In the vast theatre of life, symbiosis stands as a testament to the interconnectedness of all living beings" "Art and culture possess immense historical significance, as they provide a unique lens through which we can understand the evolution of societies, civilizations, and the human condition itself"

Replace-Text "From the mutualistic harmony of bees pollinating flowers to the exploitative nature of parasitic relationships, the web of symbiosis weaves together the tapestry of ecosystems" "From prehistoric cave paintings to modern masterpieces, art has served as a means of communication, storytelling, and self-expression, reflecting the values, beliefs, and experiences of different cultures and eras"

Replace-Text " It imparts invaluable lessons on cooperation, mutual advantage, and the delicate balance that fosters survival." ""

Replace-Text "Understanding symbiosis enhances our appreciation for the intricate beauty of nature's collaborative choreography, guiding us towards a sustainable and harmonious coexistence with the natural world" "The study of arts and culture not only enriches our understanding of history but also fosters tolerance, promotes critical thinking, and enhances our appreciation for the beauty and complexity of human creativity"

# --- Append a new empty paragraph at the very end of the document body ---
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

Write-Output "done"
